$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Dudas" (questions) table already lists 6 questions in rows 3-8.
# New team questions (7-11) are appended as rows 9-13, following the exact
# same layout/formatting as the existing rows (col A = index number,
# col B = question text, col C = empty "answered?" cell).

# Copy the formatting (styles/borders/fill) of the last existing question
# row (row 8) down onto the five new rows so the new cells pick up the same
# cell styles (centered index, left-aligned wrapped text, highlighted
# answer cell) used by every other question row.
$ws.Range("A8:H8").Copy() | Out-Null
$ws.Range("A9:H13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$questions = @(
    "Es necesario aclarar PRE y POS condiciones de cada metodo de la clase sistema",
    "Documentacion en PDF?, ante se pedia con eleccion de estructuras y fundamentos, letra no pide.",
    "Calclar max puntos del grafo, for o contador ??",
    "E.equals()nde vertice?? Es por corrdX/coordY ????",
    "NodoCritico ?"
)

$row = 9
$num = 7
foreach ($q in $questions) {
    $ws.Cells.Item($row, 1).Value = $num
    $ws.Cells.Item($row, 2).Value = $q
    $row = $row + 1
    $num = $num + 1
}

# Leave the selection where the author last clicked before saving.
$ws.Range("B14").Select()
